# Implemented Profile section validation
# Adds the OrangeHRM support e-mail address next to the "Support" row on the
# "Profile" sheet, using a left-aligned style based on the sheet's normal
# (non-hyperlink) font.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Profile")

# New cell: B3, on the same row as A3 ("Support").
$ws.Range("B3").Value = " ossupport@orangehrm.com "

# Give B3 the same base formatting (font) as the other plain data cells
# (e.g. A4 / "Changepassword"), then left-align it - this reproduces the new
# cellXfs entry (fontId=1, horizontal=left alignment) added to styles.xml.
$ws.Range("A4").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").HorizontalAlignment = -4131
